$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6-10 (group5..group9), keep only A1:A5
$ws.Range("A6:A10").EntireRow.Delete() | Out-Null

# Update A1:A5 with new random-looking group strings
$ws.Range("A1").Value = "g2844"
$ws.Range("A2").Value = "g4zUdx"
$ws.Range("A3").Value = "g"
$ws.Range("A4").Value = "gGzap"
$ws.Range("A5").Value = "g 3p"
